$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "1266001"
$ws.Range("A3").Value = "1266001"
